# "Confirmed which features I plan to implement": assign Annabelle as the
# person for the Tutorial Tab / Create Tutorials features (previously N/A),
# and leave the selection on the last cell touched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = "Annabelle"
$ws.Range("C10").Value = "Annabelle"
$ws.Range("C10").Select()
